$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update computed values in columns A and B (rows 1-32)
$ws.Range("A1").Value = -0.35360442874736009
$ws.Range("B1").Value = 0.35283020095378959
$ws.Range("A2").Value = -0.27471834282242469
$ws.Range("B2").Value = 0.27214746042431059
$ws.Range("A3").Value = -0.16919886825060004
$ws.Range("B3").Value = 0.16849953760311465
$ws.Range("A4").Value = -0.15649953777843884
$ws.Range("B4").Value = 0.15587985603348109
$ws.Range("A5").Value = -0.14987985673548021
$ws.Range("B5").Value = 0.14864875967617674
$ws.Range("A6").Value = -0.048440905681081059
$ws.Range("B6").Value = 0.048411473917091108
$ws.Range("A7").Value = -0.028411474760392963
$ws.Range("B7").Value = 0.028374061174256227
$ws.Range("A8").Value = -0.0083740620200618565
$ws.Range("B8").Value = 0.0083560936045357082
$ws.Range("A9").Value = -0.0023560943323879258
$ws.Range("B9").Value = 0.0023381434464164741
$ws.Range("A10").Value = 0.0036618558256122924
$ws.Range("B10").Value = -0.0036611376115516236
$ws.Range("A11").Value = 0.0081611368963443454
$ws.Range("B11").Value = -0.0081707274157096776
$ws.Range("A12").Value = 0.014170726688131019
$ws.Range("B12").Value = -0.014249222138232209
$ws.Range("A13").Value = -0.058561165142589466
$ws.Range("B13").Value = 0.058438062536215085
$ws.Range("A14").Value = -0.046438063315964229
$ws.Range("B14").Value = 0.046347294891166513
$ws.Range("A15").Value = -0.040347295624047597
$ws.Range("B15").Value = 0.040256632772782908
$ws.Range("A16").Value = -0.01502675052828284
$ws.Range("B16").Value = 0.015004070216859677
$ws.Range("A17").Value = -0.0090040709569620958
$ws.Range("B17").Value = 0.0089999992338700707
$ws.Range("A18").Value = -0.068111858829201566
$ws.Range("B18").Value = 0.06806628245635693
$ws.Range("A19").Value = -0.027096526772385943
$ws.Range("B19").Value = 0.027013309540879771
$ws.Range("A20").Value = -0.018013310237250835
$ws.Range("B20").Value = 0.018004259495393526
$ws.Range("A21").Value = -0.0090042601926807464
$ws.Range("B21").Value = 0.0089999993021576685
$ws.Range("A22").Value = -0.093948592883801751
$ws.Range("B22").Value = 0.093635574191473836
$ws.Range("A23").Value = -0.08463557490725826
$ws.Range("B23").Value = 0.084126920916991033
$ws.Range("A24").Value = -0.042126921924260152
$ws.Range("B24").Value = 0.041999998987352072
$ws.Range("A25").Value = -0.094207854718138861
$ws.Range("B25").Value = 0.094090628998859671
$ws.Range("A26").Value = -0.088090629723811986
$ws.Range("B26").Value = 0.087944575990810847
$ws.Range("A27").Value = -0.08194457671891664
$ws.Range("B27").Value = 0.081461636120741865
$ws.Range("A28").Value = -0.075461636863244586
$ws.Range("B28").Value = 0.075148563876126673
$ws.Range("A29").Value = -0.063148564678643382
$ws.Range("B29").Value = 0.06301261266647451
$ws.Range("A30").Value = -0.043012613542025679
$ws.Range("B30").Value = 0.042661343303768806
$ws.Range("A31").Value = -0.027661344144060962
$ws.Range("B31").Value = 0.027554955344815113
$ws.Range("A32").Value = -0.0060006435258497248
$ws.Range("B32").Value = 0.0059999992348416242

# Widen column B to match column A (15.42578125 characters)
$ws.Columns.Item(2).ColumnWidth = 14.67
